$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.228.52'
$ws.Range("E2").Value = '  -2.79%  '
$ws.Range("D3").Value = '1.549.43'
$ws.Range("E3").Value = '  -4.80%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.85'
$ws.Range("E5").Value = '  -3.41%  '
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("E7").Value = '  -5.63%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0606'
$ws.Range("E8").Value = '  -1.99%  '
$ws.Range("E9").Value = '  -3.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '17.77'
$ws.Range("E10").Value = '  -3.71%  '
$ws.Range("E11").Value = '  -1.16%  '
$ws.Range("E12").Value = '  -4.84%  '
$ws.Range("D13").Value = '1.546.34'
$ws.Range("E13").Value = '  -5.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.98'
$ws.Range("E14").Value = '  -4.87%  '
$ws.Range("E15").Value = '  -4.68%  '
$ws.Range("D16").Value = '25.184.54'
$ws.Range("E16").Value = '  -3.01%  '
$ws.Range("E17").Value = '  -4.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '58.46'
$ws.Range("E18").Value = '  -4.69%  '
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '185.70'
$ws.Range("E20").Value = '  -3.37%  '
$ws.Range("E21").Value = '  -3.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.23'
$ws.Range("E22").Value = '  -3.69%  '
$ws.Range("E23").Value = '  -4.28%  '
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("E25").Value = '  -4.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '139.61'
$ws.Range("E26").Value = '  -2.91%  '
$ws.Range("E27").Value = '  -4.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '14.76'
$ws.Range("E28").Value = '  -2.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.38'
$ws.Range("E29").Value = '  -5.26%  '
$ws.Range("E30").Value = '  -6.76%  '
$ws.Range("E31").Value = '  -4.59%  '
$ws.Range("E32").Value = '  -3.91%  '
$ws.Range("E33").Value = '  -4.81%  '
$ws.Range("E34").Value = '  -3.40%  '
$ws.Range("E35").Value = '  -4.14%  '
$ws.Range("D36").Value = '1.084.56'
$ws.Range("E36").Value = '  -3.69%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("E38").Value = '  -2.70%  '
$ws.Range("E39").Value = '  -5.46%  '
$ws.Range("E40").Value = '  -7.47%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.759'
$ws.Range("E41").Value = '  -10.97%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.799'
$ws.Range("E42").Value = '  +4.64%  '
$ws.Range("E43").Value = '  -5.73%  '
$ws.Range("E44").Value = '  -2.59%  '
$ws.Range("D45").Value = '1.680.60'
$ws.Range("E45").Value = '  -4.83%  '
$ws.Range("E46").Value = '  +13.75%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.44'
$ws.Range("E47").Value = '  -2.12%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '52.20'
$ws.Range("E48").Value = '  -4.03%  '
$ws.Range("E49").Value = '  -5.95%  '
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("E51").Value = '  -2.05%  '
